$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.293.17'
$ws.Range("E2").Value = '  -1.60%  '
$ws.Range("D3").Value = '2.176.64'
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("E6").Value = '  -1.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '70.37'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.58%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.42'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0928'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("E12").Value = '  -5.91%  '
$ws.Range("E13").Value = '  -4.54%  '
$ws.Range("E14").Value = '  -1.89%  '
$ws.Range("D15").Value = '2.498.79'
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("E17").Value = '  -4.19%  '
$ws.Range("D18").Value = '2.176.03'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '41.100.03'
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("E20").Value = '  -6.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("E22").Value = '  -2.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '226.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.13%  '
$ws.Range("E25").Value = '  -5.70%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.35%  '
$ws.Range("E29").Value = '  -2.38%  '
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0771'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.18'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.82%  '
$ws.Range("E36").Value = '  -2.68%  '
$ws.Range("E37").Value = '  -6.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.74%  '
$ws.Range("E39").Value = '  -4.63%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.98'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.46%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '59.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.35%  '
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("E45").Value = '  -2.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.20%  '
$ws.Range("E48").Value = '  -1.72%  '
$ws.Range("E49").Value = '  -1.68%  '
$ws.Range("E50").Value = '  -6.85%  '
